$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colG = $ws.Range("G1:G259")
$colG.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System") | Out-Null
